$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 74, pushing the existing "Segunda" row (old row 74)
# down to row 75. The new (blank) row becomes row 74.
$ws.Rows.Item(74).Insert()

# Row 73 now reflects the latest weekly report for the "Primera" quality grade:
# new date, new min/max/avg/kg prices (volume, unit and quality stay the same).
$ws.Range("D73").Value = 45121
$ws.Range("N73").Value = 3000
$ws.Range("O73").Value = 3500
$ws.Range("P73").Value = 3200
$ws.Range("S73").Value = 320

# New row 74 carries what used to be row 73's data (previous week's "Primera"
# quality entry), now recorded as its own historical row.
$ws.Range("A74").Value = 1
$ws.Range("B74").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C74").Value = "Arica y Parinacota"
$ws.Range("D74").Value = 45072
$ws.Range("E74").Value = 15
$ws.Range("F74").Value = "Fruta"
$ws.Range("G74").Value = 100108
$ws.Range("H74").Value = "Tropicales y subtropicales"
$ws.Range("I74").Value = 100108001
$ws.Range("J74").Value = "Guayaba"
$ws.Range("K74").Value = "Sin especificar"
$ws.Range("L74").Value = "Primera"
$ws.Range("M74").Value = 250
$ws.Range("N74").Value = 5000
$ws.Range("O74").Value = 6000
$ws.Range("P74").Value = 5600
$ws.Range("Q74").Value = "$/caja 10 kilos"
$ws.Range("R74").Value = "Región de Arica y Parinacota"
$ws.Range("S74").Value = 560
$ws.Range("T74").Value = 10
